# Daily attendance processing - 2026-01-22 08:44:19
# Swap the order of the first two comma-separated "Recorded By" entries
# in column G for the specified rows (leaving any trailing entries, e.g.
# the extra "system" token, untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,7,8,28,29,30,31,32,33,34,54,55,56,57,58,59,60,80,81,82,87,106,107,108,113,132,133,134,139)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2
    if ($val -ne $null) {
        $parts = $val -split ', '
        if ($parts.Length -ge 2) {
            $tmp = $parts[0]
            $parts[0] = $parts[1]
            $parts[1] = $tmp
            $cell.Value = [string]::Join(', ', $parts)
        }
    }
}
